$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # quality_comparison
$ws2 = $wb.Worksheets.Item(2)   # computational_comparison

# --- 1. Rename "fedcore" -> "approach" in header rows ---
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# --- 2. Apply new border styles to the merged-header spacer cells ---
# Build the two new styles once (on sheet1 C1/D1), then replicate them
# via copy/paste-format so the style table stays clean (no stray xfs).

$c1 = $ws1.Range("C1")
$c1.Style = "Normal"
$c1.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$c1.Borders.Item(9).LineStyle = 1   # xlEdgeBottom

$d1 = $ws1.Range("D1")
$d1.Style = "Normal"
$d1.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$d1.Borders.Item(10).LineStyle = 1  # xlEdgeRight
$d1.Borders.Item(9).LineStyle = 1   # xlEdgeBottom

# Replicate the same two styles onto sheet2's C1/D1 and F1/G1
$c1.Copy()
$ws2.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$ws2.Range("F1").PasteSpecial(-4122)

$d1.Copy()
$ws2.Range("D1").PasteSpecial(-4122)
$ws2.Range("G1").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- 3. Remove the stray empty inline-string cell G5 on sheet2 ---
$ws2.Range("G5").ClearContents()
